$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Addr,
        [string]$Text
    )
    $rng = $ws.Range($Addr)
    # Force text format so numeric-looking strings (e.g. "384.15") are not
    # silently coerced into floating point numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# Row 2 - Bitcoin
Set-CellText "D2" "51.367.98"
Set-CellText "E2" "  +0.41%  "

# Row 3 - Ethereum
Set-CellText "D3" "3.040.43"
Set-CellText "E3" "  +2.68%  "

# Row 4 - TetherUSD
Set-CellText "D4" "0.999"
Set-CellText "E4" "  -0.03%  "

# Row 5 - BNB
Set-CellText "D5" "384.15"
Set-CellText "E5" "  +1.05%  "

# Row 6 - Solana
Set-CellText "D6" "102.64"
Set-CellText "E6" "  +0.29%  "

# Row 7
Set-CellText "D7" "0.542"
Set-CellText "E7" "  -0.42%  "

# Row 8
Set-CellText "E8" "  -0.01%  "

# Row 9
Set-CellText "E9" "  -1.14%  "

# Row 10
Set-CellText "D10" "36.73"
Set-CellText "E10" "  +0.70%  "

# Row 11
Set-CellText "D11" "0.137"
Set-CellText "E11" "  +0.13%  "

# Row 12
Set-CellText "D12" "0.0861"
Set-CellText "E12" "  +0.17%  "

# Row 13
Set-CellText "D13" "3.518.40"
Set-CellText "E13" "  +2.68%  "

# Row 14
Set-CellText "D14" "18.55"
Set-CellText "E14" "  +1.53%  "

# Row 15
Set-CellText "D15" "7.73"
Set-CellText "E15" "  -0.55%  "

# Row 16
Set-CellText "D16" "3.034.62"
Set-CellText "E16" "  +2.18%  "

# Row 17
Set-CellText "D17" "0.966"
Set-CellText "E17" "  -2.50%  "

# Row 18
Set-CellText "D18" "10.55"
Set-CellText "E18" "  -5.65%  "

# Row 19
Set-CellText "D19" "51.411.12"
Set-CellText "E19" "  +0.43%  "

# Row 20
Set-CellText "D20" "3.13"
Set-CellText "E20" "  -0.88%  "

# Row 21
Set-CellText "D21" "12.39"
Set-CellText "E21" "  -1.00%  "

# Row 22
Set-CellText "D22" "0.0₃0963"
Set-CellText "E22" "  +0.35%  "

# Row 23
Set-CellText "D23" "70.06"
Set-CellText "E23" "  -0.03%  "

# Row 24
Set-CellText "D24" "266.87"
Set-CellText "E24" "  +0.15%  "

# Row 25
Set-CellText "E25" "  -1.55%  "

# Row 26
Set-CellText "D26" "8.16"
Set-CellText "E26" "  +4.54%  "

# Row 27
Set-CellText "D27" "26.83"
Set-CellText "E27" "  +3.45%  "

# Row 28
Set-CellText "D28" "0.169"
Set-CellText "E28" "  +2.70%  "

# Row 29
Set-CellText "D29" "7.24"
Set-CellText "E29" "  -2.44%  "

# Row 30
Set-CellText "E30" "  +0.10%  "

# Row 31
Set-CellText "E31" "  -1.85%  "

# Row 32
Set-CellText "D32" "10.26"
Set-CellText "E32" "  -0.46%  "

# Row 33
Set-CellText "D33" "34.89"
Set-CellText "E33" "  +0.86%  "

# Row 34
Set-CellText "E34" "  +2.91%  "

# Row 35
Set-CellText "D35" "50.25"
Set-CellText "E35" "  -1.85%  "

# Row 36
Set-CellText "E36" "  +1.88%  "

# Row 37
Set-CellText "E37" "  -0.15%  "

# Row 38
Set-CellText "D38" "3.34"
Set-CellText "E38" "  +3.12%  "

# Row 39
Set-CellText "D39" "0.289"
Set-CellText "E39" "  +6.94%  "

# Row 40
Set-CellText "D40" "16.96"
Set-CellText "E40" "  +2.14%  "

# Row 41
Set-CellText "E41" "  +1.30%  "

# Row 42
Set-CellText "E42" "  -0.79%  "

# Row 43
Set-CellText "E43" "  -0.01%  "

# Row 44
Set-CellText "D44" "124.56"
Set-CellText "E44" "  -0.16%  "

# Row 45
Set-CellText "E45" "  +4.02%  "

# Row 46
Set-CellText "D46" "21.77"
Set-CellText "E46" "  +1.59%  "

# Row 47 - now WEMIXToken (was ApeXProtocol)
Set-CellText "B47" "WEMIXToken"
Set-CellText "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText "D47" "2.08"
Set-CellText "E47" "  +2.76%  "

# Row 48 - now ApeXProtocol (was WEMIXToken)
Set-CellText "B48" "ApeXProtocol"
Set-CellText "C48" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-CellText "D48" "2.42"
Set-CellText "E48" "  +1.58%  "

# Row 49 - Maker
Set-CellText "D49" "2.024.32"
Set-CellText "E49" "  -0.50%  "

# Row 50 - RocketPoolETH
Set-CellText "E50" "  +2.54%  "

# Row 51 - now BEAM (was Algorand)
Set-CellText "B51" "BEAM"
Set-CellText "C51" "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-CellText "D51" "0.0317"
Set-CellText "E51" "  -3.32%  "
